$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Stack ops & interrupts" -> "Stack ops"; drop the BRK/RTI line ---
$ws.Range("A6").Value = "Stack ops"
$ws.Range("B6").Value = "- PHA, PLA`n- PHP, PLP"
$ws.Rows.Item(6).RowHeight = 24.25

# --- New row 12: Google Test / URL (hyperlink) / Done ---
$ws.Range("A12").Value = "Google Test"
$ws.Range("C12").Value = "Done"

$linkCell = $ws.Range("B12")
$linkCell.Value = "https://github.com/google/googletest.git"
$ws.Hyperlinks.Add($linkCell, "https://github.com/google/googletest.git", "", "", "https://github.com/google/googletest.git")
$font = $linkCell.Font
$font.Color = 16711680
$font.Underline = $false
$linkCell.WrapText = $true
$linkCell.VerticalAlignment = -4160

# --- Selection moves to B9 ---
$null = $ws.Range("B9").Select()

Write-Output "done"
